# Update "想去人数" (want-to-go count) figures for a handful of events.
# These numbers are refreshed by the site's automated gh-pages data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5589   # was 5585
$wsExhibit.Range("F6").Value = 833    # was 832
$wsExhibit.Range("F7").Value = 51     # was 50
$wsExhibit.Range("F8").Value = 373    # was 372
$wsExhibit.Range("F10").Value = 4     # was 3

# --- Sheet "全部类型" (same events, different row offsets) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5589   # was 5585
$wsAll.Range("F6").Value = 833    # was 832
$wsAll.Range("F7").Value = 51     # was 50
$wsAll.Range("F9").Value = 373    # was 372
$wsAll.Range("F11").Value = 4     # was 3
